# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" sheet (1st sheet) and the "全部类型" sheet (4th sheet), to
# reflect newly scraped counts (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -----------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2, 6).Value  = 188
$ws1.Cells.Item(5, 6).Value  = 959
$ws1.Cells.Item(6, 6).Value  = 5221
$ws1.Cells.Item(7, 6).Value  = 438
$ws1.Cells.Item(8, 6).Value  = 619
$ws1.Cells.Item(9, 6).Value  = 905
$ws1.Cells.Item(13, 6).Value = 563
$ws1.Cells.Item(14, 6).Value = 10
$ws1.Cells.Item(17, 6).Value = 1745
$ws1.Cells.Item(18, 6).Value = 1448
$ws1.Cells.Item(19, 6).Value = 816
$ws1.Cells.Item(21, 6).Value = 186
$ws1.Cells.Item(23, 6).Value = 504
$ws1.Cells.Item(28, 6).Value = 2510
$ws1.Cells.Item(30, 6).Value = 96
$ws1.Cells.Item(31, 6).Value = 49
$ws1.Cells.Item(33, 6).Value = 22
$ws1.Cells.Item(34, 6).Value = 247
$ws1.Cells.Item(39, 6).Value = 272
$ws1.Cells.Item(41, 6).Value = 80

# --- Sheet 4: 全部类型 --------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(3, 6).Value  = 188
$ws4.Cells.Item(5, 6).Value  = 959
$ws4.Cells.Item(7, 6).Value  = 5221
$ws4.Cells.Item(8, 6).Value  = 438
$ws4.Cells.Item(9, 6).Value  = 619
$ws4.Cells.Item(12, 6).Value = 905
$ws4.Cells.Item(17, 6).Value = 563
$ws4.Cells.Item(18, 6).Value = 10
$ws4.Cells.Item(22, 6).Value = 1745
$ws4.Cells.Item(23, 6).Value = 1448
$ws4.Cells.Item(24, 6).Value = 816
$ws4.Cells.Item(26, 6).Value = 186
$ws4.Cells.Item(29, 6).Value = 504
$ws4.Cells.Item(33, 6).Value = 2510
$ws4.Cells.Item(35, 6).Value = 96
$ws4.Cells.Item(37, 6).Value = 22
$ws4.Cells.Item(38, 6).Value = 247
$ws4.Cells.Item(42, 6).Value = 272
$ws4.Cells.Item(44, 6).Value = 80

$wb.Save()
